# Weekly data refresh: insert the newest week's price record for
# "Haba" (Mercado Mayorista Lo Valledor de Santiago) ahead of the
# existing history, pushing every later row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 429; Excel shifts rows 429:501 down to 430:502
# and grows the used range (dimension) accordingly.
$ws.Rows(429).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A429").Value = 6
$ws.Range("B429").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C429").Value = "Metropolitana"
$ws.Range("D429").Value = 45258
$ws.Range("E429").Value = 13
$ws.Range("F429").Value = 100112026
$ws.Range("G429").Value = "Haba"
$ws.Range("H429").Value = "Sin especificar"
$ws.Range("I429").Value = "Primera"
$ws.Range("J429").Value = 770
$ws.Range("K429").Value = 7000
$ws.Range("L429").Value = 8000
$ws.Range("M429").Value = 7416
$ws.Range("N429").Value = "`$/saco 25 kilos"
$ws.Range("O429").Value = "Región del Maule"
$ws.Range("P429").Value = 297
$ws.Range("Q429").Value = 25
$ws.Range("R429").Value = "Hortaliza"
